$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.933.39"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "1.643.25"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.72"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0874"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").Value = "1.875.95"
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("D13").Value = "1.637.15"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.572"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.15%  "
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "27.912.52"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.75%  "
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("E24").Value = "  +2.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("E32").Value = "  +2.08%  "
$ws.Range("D33").Value = "1.427.25"
$ws.Range("E33").Value = "  -2.63%  "
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("E37").Value = "  +1.89%  "
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.926"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.560"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "68.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.37%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.97%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "1.784.76"
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "89.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("E51").Value = "  +0.59%  "
